$p = $ppt.ActivePresentation

# --- Slide 13: "Create a New Project" ---
# "...odify a previously created project using again the project tool:"
# -> "...odify a previously created project by using the project tool again:"
$s13 = $p.Slides.Item(13)
$sh13 = $s13.Shapes.Item(2)
$tr13 = $sh13.TextFrame.TextRange
$run13 = $tr13.Characters(196, 64)
$run13.Text = "odify a previously created project by using the project tool again:"

# --- Slide 16: "Install Your Project" ---
$s16 = $p.Slides.Item(16)
$sh16 = $s16.Shapes.Item(2)
$tr16 = $sh16.TextFrame.TextRange

# Merge "-based projects are build this way. " + "The " into a single run
$run16a = $tr16.Characters(129, 40)
$run16a.Text = "-based projects are build this way. The "

# Merge "summarizes " + "these steps." into a single run
$run16b = $tr16.Characters(199, 23)
$run16b.Text = "summarizes these steps."

# "Build and install the (yet empty) project:" -> "Build and install the (currently empty) project:"
$run16c = $tr16.Characters(224, 42)
$run16c.Text = "Build and install the (currently empty) project:"
